$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns are treated as plain text so
# numeric-looking strings (e.g. "52.009.32") are not auto-coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$updates = @{
    2  = @{ D = "52.009.32";  E = "  -1.36%  " }
    3  = @{ D = "2.809.47";   E = "  -1.39%  " }
    4  = @{ D = "1.00";       E = "  +0.32%  " }
    5  = @{ D = "358.33";     E = "  -0.54%  " }
    6  = @{ D = "110.37";     E = "  -5.66%  " }
    7  = @{ D = "0.562";      E = "  +1.95%  " }
    8  = @{ D = "1.00";       E = "  +0.18%  " }
    9  = @{ D = "0.595";      E = "  -1.96%  " }
    10 = @{ D = "40.28";      E = "  -6.58%  " }
    11 = @{ D = "0.0855";     E = "  -1.56%  " }
    12 = @{ D = "0.133";      E = "  +1.43%  " }
    13 = @{ D = "19.70";      E = "  -2.29%  " }
    14 = @{ D = "7.72";       E = "  -2.77%  " }
    15 = @{ D = "3.262.44";   E = "  -0.60%  " }
    16 = @{ D = "2.877.51";   E = "  +1.40%  " }
    17 = @{ D = "0.914";      E = "  +0.78%  " }
    18 = @{ D = "52.011.65";  E = "  -1.16%  " }
    19 = @{ D = "7.46";       E = "  +1.93%  " }
    20 = @{ D = "3.12";       E = "  -2.57%  " }
    21 = @{ D = "13.22";      E = "  -3.78%  " }
    22 = @{ D = "0.0₃0989";   E = "  -0.29%  " }
    23 = @{ D = "272.08";     E = "  -0.52%  " }
    24 = @{ D = "69.92";      E = "  -1.21%  " }
    25 = @{ D = "2.80";       E = "  -1.80%  " }
    26 = @{ D = "26.65";      E = "  -2.68%  " }
    27 = @{ E = "  +0.01%  " }
    28 = @{ D = "10.21";      E = "  -1.70%  " }
    29 = @{ E = "  -1.18%  " }
    30 = @{ D = "0.141";      E = "  -0.55%  " }
    31 = @{ D = "0.0473";     E = "  +2.70%  " }
    32 = @{ D = "52.08";      E = "  +1.39%  " }
    33 = @{ D = "34.11";      E = "  -1.97%  " }
    34 = @{ D = "5.79";       E = "  -1.41%  " }
    35 = @{ D = "5.48";       E = "  +10.31%  " }
    36 = @{ D = "0.0842";     E = "  -0.07%  " }
    37 = @{ D = "1.00";       E = "  +0.22%  " }
    38 = @{ D = "3.19";       E = "  -2.79%  " }
    39 = @{ D = "2.01";       E = "  -5.84%  " }
    40 = @{ D = "17.97";      E = "  -5.15%  " }
    41 = @{ D = "0.116";      E = "  -0.26%  " }
    42 = @{ D = "2.52";       E = "  -4.79%  " }
    43 = @{ D = "125.81";     E = "  -0.43%  " }
    44 = @{ E = "  -1.15%  " }
    45 = @{ D = "22.47";      E = "  -5.49%  " }
    46 = @{ D = "2.062.30";   E = "  -0.74%  " }
    47 = @{ D = "3.26";       E = "  -3.98%  " }
    48 = @{ D = "2.33";       E = "  +1.20%  " }
    49 = @{ D = "5.82";       E = "  +2.41%  " }
    50 = @{ D = "0.942";      E = "  -2.80%  " }
    51 = @{ D = "9.05";       E = "  +0.00%  " }
}

foreach ($rowNum in $updates.Keys) {
    $rowData = $updates[$rowNum]
    if ($rowData.ContainsKey("D")) {
        $ws.Range("D$rowNum").Value = $rowData["D"]
    }
    $ws.Range("E$rowNum").Value = $rowData["E"]
}
